$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Holly added S.GISH to the harvester column (B) for all data rows.
$ws.Range("B2:B18").Value = "S.GISH"

# Reflect the slightly adjusted column B width from the autofit/edit.
$ws.Columns("B").ColumnWidth = 8

# Leave the selection on column B, matching the editor's last action.
$ws.Range("B:B").Select() | Out-Null
